$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 93, pushing the existing rows 93-116 down to 94-117,
# carrying along their formatting (matches the diff: row count grows to 117).
$ws.Rows("93:93").Insert()

# Populate the newly inserted row 93 with the new record.
$ws.Range("A93").Value = 1
$ws.Range("B93").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C93").Value = "Arica y Parinacota"
$ws.Range("D93").Value = 44943
$ws.Range("D93").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E93").Value = 15
$ws.Range("F93").Value = 100112038
$ws.Range("G93").Value = "Cebollín baby"
$ws.Range("H93").Value = "Sin especificar"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 200
$ws.Range("K93").Value = 5000
$ws.Range("L93").Value = 6000
$ws.Range("M93").Value = 5500
$ws.Range("N93").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O93").Value = "Región de Arica y Parinacota"
$ws.Range("P93").Value = 2750
$ws.Range("Q93").Value = 2
$ws.Range("R93").Value = "Hortaliza"
